$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @(
    "Menu_Apperance_ThemeAuto.dds",
    "Menu_Apperance_ThemeManual.dds",
    "Menu_Info.dds",
    "Menu_Device.dds",
    "Menu_Misc.dds"
)

$startRow = 29
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $row
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

$ws.Range("B33").Select() | Out-Null
